$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the "Desarquivamentos Pendentes" sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Keep "PAINEIS DARQ" as the active/selected sheet
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
